$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-20 08:19:09"
$wsZhCn.Range("G5").Value = "2016-01-20 08:19:56"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-20 08:19:20"
$wsDeDe.Range("G5").Value = "2016-01-20 08:20:16"
